$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.472.98"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "'1.555.28"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("D5").Value = "'210.65"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").Value = "'0.484"
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "'24.23"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").Value = "'1.777.86"
$ws.Range("D13").Value = "'1.564.55"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").Value = "'28.472.82"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").Value = "'61.17"
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").Value = "'229.28"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").Value = "'0.0₃0671"
$ws.Range("E20").Value = "  -2.16%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "'3.91"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("E23").Value = "  -1.74%  "
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").Value = "'151.05"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  -1.35%  "
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("E29").Value = "  -2.37%  "
$ws.Range("D30").Value = "'0.0463"
$ws.Range("E30").Value = "  -3.70%  "
$ws.Range("E31").Value = "  -3.12%  "
$ws.Range("D32").Value = "'3.17"
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("D33").Value = "'1.394.58"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'3.00"
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("E35").Value = "  -2.84%  "
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("E37").Value = "  -2.63%  "
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("E44").Value = "  +1.78%  "
$ws.Range("E45").Value = "  +3.96%  "
$ws.Range("D46").Value = "'5.32"
$ws.Range("E46").Value = "  -1.65%  "
$ws.Range("D47").Value = "'1.690.75"
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("D48").Value = "'0.867"
$ws.Range("E48").Value = "  -6.46%  "
$ws.Range("E49").Value = "  +6.46%  "
$ws.Range("D50").Value = "'85.40"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("E51").Value = "  +0.02%  "
